# "pontos notáveis - incremento na tabela de ranking"
# Scale the "particip" (E) and "taxa_sucesso" (F) ranking columns by 100
# (rows 2-74 of the summary table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 74; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # column E: particip
    $fCell = $ws.Cells.Item($r, 6)   # column F: taxa_sucesso

    $eCell.Value = $eCell.Value2 * 100
    $fCell.Value = $fCell.Value2 * 100
}
